$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.261.51"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "2.588.03"
$ws.Range("E3").Value = "  +1.61%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("E8").Value = "  +2.08%  "
$ws.Range("E9").Value = "  +4.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.67"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.358"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.85%  "
$ws.Range("D14").Value = "3.053.67"
$ws.Range("E14").Value = "  +1.76%  "
$ws.Range("D15").Value = "63.206.13"
$ws.Range("E15").Value = "  +0.27%  "
$ws.Range("E16").Value = "  +4.03%  "
$ws.Range("D17").Value = "2.576.62"
$ws.Range("E17").Value = "  +1.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.43"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "344.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.72%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.84%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("E23").Value = "  -3.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.42%  "
$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.171"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.54%  "
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "2.665.52"
$ws.Range("E26").Value = "  -0.31%  "
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +12.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.45%  "
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("B31").Value = "SuiNetwork"
$ws.Range("C31").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.49"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.61%  "
$ws.Range("D33").Value = "0.0₃0832"
$ws.Range("E33").Value = "  +1.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "467.04"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +12.82%  "
$ws.Range("E35").Value = "  +4.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "176.92"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.55%  "
$ws.Range("E37").Value = "  +2.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.32"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.66"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.18%  "
$ws.Range("E41").Value = "  +0.66%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "152.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.86"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.58%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0554"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.618"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0979"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0242"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.77"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.39"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.72%  "
